$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New_Mapping")

# --- Chapter 4 & 5 "About" rows (19-20) are now done; clear their outstanding/
#     needs-attention marker (red highlight in column A) and the highlight fill
#     on B:D, while keeping their text/values and borders intact. ---
$ws.Range("A19:A20").Clear()
$ws.Range("B19:D20").Interior.Pattern = -4142

# --- Append the new Ch4 & Ch5 supporting-document rows to the bottom of the
#     table (rows 123-139). Column B ids, C form, D description - matching the
#     existing table layout (columns inherit the sheet's default row styles). ---
$newRows = @(
    @("125", "MainForm",   "Put Results into Action"),
    @("126", "MainForm",   "Storytelling on Project Implementation"),
    @("127", "MainForm",   "Define the Problems to Solve"),
    @("128", "MainForm",   "Workshop 3- Results to Action"),
    @("129", "MainForm",   "Warm Up Activity"),
    @("130", "MainForm",   "Create Vision Statements"),
    @("131", "MainForm",   "Brainstorm Actions"),
    @("132", "MainForm",   "Evaluate Actions"),
    @("133", "MainForm",   "Strategy Planning"),
    @("134", "MainForm",   "Workshop Wrap-Up"),
    @("135", "MainForm",   "Ch 4 Reflection"),
    @("136", "OutputForm", "Key Takeaways from Strategizing"),
    @("138", "MainForm",   "Document Your Project"),
    @("139", "MainForm",   "Next Steps for Outreach and Engagement"),
    @("140", "MainForm",   "Next Steps for Implementation"),
    @("141", "MainForm",   "Use Monitoring to Support Future Action"),
    @("142", "MainForm",   "Ch 5 Reflection")
)

$startRow = 123
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
}

# --- Reflect the user's final on-screen position: scrolled down and landed on
#     the last cell they typed into. ---
$ws.Activate()
$ws.Range("D123").Select()
